# "Generate Report for Handback"
#
# The localization-status report is refreshed once the de-de/zh-cn
# handback round-trips complete:
#   - Overview sheet: the status caption moves from "Ready for handoff"
#     to "Handed back: in sync with en-US" (shared by both language
#     status columns E/F, which also grow wider to fit it).
#   - zh-cn sheet: Latest Target File (I2) now links to the source .md,
#     Latest Handback File (J2) records the generated .zh-cn.xlf.
#   - de-de sheet: same shape, plus a fresh Latest Handback DateTime
#     (K2) since that handback finished after the zh-cn one.

$wb = $excel.ActiveWorkbook

$mdDisplay = "b9b8b0b6-6e53-45b8-bf1b-cdc454d0f50a.md"
$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/92d2561af4499624546d3a472fa293a23b88d361/e2e/b9b8b0b6-6e53-45b8-bf1b-cdc454d0f50a.md"

$newStatus = "Handed back: in sync with en-US"

# --- zh-cn sheet ------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Columns.Item(3).ColumnWidth = 29.166666666666668
$ws2.Columns.Item(9).ColumnWidth = 39.166666666666664
$ws2.Columns.Item(10).ColumnWidth = 39.166666666666664

# Status (C2) shares its text with the Overview sheet's per-language
# status cells - updating it here is what actually flips that shared
# string everywhere it is used.
$ws2.Range("C2").Value = $newStatus

$ws2.Range("I2").Value = $mdDisplay
$ws2.Range("J2").Value = "b9b8b0b6-6e53-45b8-bf1b-cdc454d0f50a.010bfedd6e8ea978eb1e60f9dfd1536063fdd740.zh-cn.xlf"
$ws2.Range("K2").Value = "2016-08-25 20:59:31"

$ws2.Hyperlinks.Add($ws2.Range("I2"), $mdUrl, "", "", $mdDisplay)

# --- de-de sheet --------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Columns.Item(3).ColumnWidth = 29.166666666666668
$ws3.Columns.Item(9).ColumnWidth = 39.166666666666664
$ws3.Columns.Item(10).ColumnWidth = 39.166666666666664

$ws3.Range("C2").Value = $newStatus

$ws3.Range("I2").Value = $mdDisplay
$ws3.Range("J2").Value = "b9b8b0b6-6e53-45b8-bf1b-cdc454d0f50a.010bfedd6e8ea978eb1e60f9dfd1536063fdd740.de-de.xlf"
$ws3.Range("K2").Value = "2016-08-25 20:59:38"

$ws3.Hyperlinks.Add($ws3.Range("I2"), $mdUrl, "", "", $mdDisplay)

# --- Overview sheet -------------------------------------------------
# E2/F2 mirror the same "Handed back: ..." status text (same shared
# string as zh-cn!C2 / de-de!C2 above), which is also why they need to
# grow to the same width.
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("E2").Value = $newStatus
$ws1.Range("F2").Value = $newStatus

$ws1.Columns.Item(5).ColumnWidth = 29.166666666666668
$ws1.Columns.Item(6).ColumnWidth = 29.166666666666668
